$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Alecxperdu, DTP, DramaPanda, Maitredoudou, FooD_Flo",
    "Chreet, Joshua, Gririsu, Hugo, Thib, Julien, Chreet_S5, Mickey_S8, Sunka",
    "Saucisson, The soulless, Ranzyo_S5, Killian, Gwendal_S7, Julien_S7, Waikato, Mickey, Pilouche",
    "Bilel_S4, Bylost_S4, And, Xiao, Nikolas_S5, Gririsu_S5, Thib_S5, Nath_S6, Grenzo, Clovis, Thynael, Cosmos, Piiskoo, Chopa, Cha, Z4tix_S9",
    "Kamiga, Nelson, Mark, Alecxperdu_S2, Gobou, Polo, Skowa, Mark_S2",
    "Z4tix, Nikolas, Gwendal, Clem fair play, Timeo, Amaury, Schweppes, Malkovan",
    "Black Lolo, Eneko, FooD_Flo_S5, Enze, Thyx, Doggydog_S6, Enze_S7, Tim, Luc",
    "Nelson_S4, Alex., Skowa_S4, Raphi, Xiao_S5, Manta, Dahmi1 Arti, Gobou_S6, Lilian, Mehdiiii, Corentin, Jilink, Alan_S9",
    "Krak, Nwog, Ethan, Jehovah, Nicovid, Alexadventure, Step, Feiik, Tiff, Ladoly, Natoxe, Bilal, Flau_S4, Didine, Deku, Sayo, Line_S5, Bilel_S5, Clem fair play_S5, TimeoGnc",
    "Kamiga_S5, GiulfeuYT_S6, Alan_S7, Sy_boulette, Reiko, Ju, Mielle",
    "Saucisson_S5, Malkovan_S5, Schweppes_S6, Ethan_S6, Xori, Crypto, Yatho, Vah Balress, Gigi, Quentin, Mistimat, Spider, Enze_S8, Alex, Guigui_S9",
    "Cyrf, Sneus, Malkovan_S3, Bylost, Line, HiYoucef, Julien_S5, Timeo_S5, Black Lolo_S5, Amaury_S5, Kwinn",
    "Lili, Samson, SorciShoot, Louan, Xo, Napoleon, Cyrf_S4, Espoir Perdu, And_S5, Major Chris, Nikolas_S6, Nemocca, Aurel, Theo, Chatoon, Clement, Mizuki, Armand, Manu",
    "Yuu, Nounours, Saminette, Cyrf_S6, Blgham, Blacks Star, Kamiga_S9",
    "Bilel, Joshua_S3, GiulfeuYT, Doggydog, Alan, Mark_S5, Maitredoudou_S5, Nicovid_S6",
    "Flau, Nath, Ethan_S3, Ranzyo, Guigui, Angel, Mtking, Dragon"
)

# Update rows 2..17 (A = 1..16, B = new text) with the new data
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove the now-unused rows 18..26 (old rows 18-26 no longer exist)
$ws.Rows("18:26").Delete()
